$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Ranking")
$ws.Range("H7").Value = 2.360530197348857
$ws.Range("I7").Value = 1.507753183552359
$ws.Range("H8").Value = 2.742151682354039
$ws.Range("I8").Value = 1.791334535278701
$ws.Range("H10").Value = 3.061910308268288
$ws.Range("I10").Value = 1.904383464632066

$ws = $wb.Worksheets.Item("Matriz_Pvalores")
$ws.Range("G2").Value = 0.06546638157961904
$ws.Range("H2").Value = 0.2296290742447789
$ws.Range("J2").Value = 0.4370194909566052
$ws.Range("G3").Value = 0.3196842051293389
$ws.Range("H3").Value = 0.3611971997491434
$ws.Range("J3").Value = 0.8714115499384436
$ws.Range("G4").Value = 0.3010811264265782
$ws.Range("H4").Value = 0.3308286524259445
$ws.Range("J4").Value = 0.53251175688164
$ws.Range("G5").Value = 0.1927496505398922
$ws.Range("H5").Value = 0.3733633158884433
$ws.Range("J5").Value = 0.5675821404899604
$ws.Range("G6").Value = 0.1151776154665982
$ws.Range("H6").Value = 0.3178144093142117
$ws.Range("J6").Value = 0.5769710716092553
$ws.Range("B7").Value = 0.06546638157961904
$ws.Range("C7").Value = 0.3196842051293389
$ws.Range("D7").Value = 0.3010811264265782
$ws.Range("E7").Value = 0.1927496505398922
$ws.Range("F7").Value = 0.1151776154665982
$ws.Range("H7").Value = 0.2890165462237508
$ws.Range("I7").Value = 0.07858805481375941
$ws.Range("J7").Value = 0.117579344021491
$ws.Range("B8").Value = 0.2296290742447789
$ws.Range("C8").Value = 0.3611971997491434
$ws.Range("D8").Value = 0.3308286524259445
$ws.Range("E8").Value = 0.3733633158884433
$ws.Range("F8").Value = 0.3178144093142117
$ws.Range("G8").Value = 0.2890165462237508
$ws.Range("I8").Value = 0.5213751693337416
$ws.Range("J8").Value = 0.208987530822311
$ws.Range("G9").Value = 0.07858805481375941
$ws.Range("H9").Value = 0.5213751693337416
$ws.Range("J9").Value = 0.7166267432886311
$ws.Range("B10").Value = 0.4370194909566052
$ws.Range("C10").Value = 0.8714115499384436
$ws.Range("D10").Value = 0.53251175688164
$ws.Range("E10").Value = 0.5675821404899604
$ws.Range("F10").Value = 0.5769710716092553
$ws.Range("G10").Value = 0.117579344021491
$ws.Range("H10").Value = 0.208987530822311
$ws.Range("I10").Value = 0.7166267432886311

$ws = $wb.Worksheets.Item("Matriz_DM_Original")
$ws.Range("G2").Value = 2.518427412732279
$ws.Range("H2").Value = 1.416312553198512
$ws.Range("J2").Value = 0.8625711600869045
$ws.Range("G3").Value = 1.13524266567742
$ws.Range("H3").Value = 1.030039409186893
$ws.Range("J3").Value = 0.1725125760959997
$ws.Range("G4").Value = 1.186498229823468
$ws.Range("H4").Value = 1.105839399979577
$ws.Range("J4").Value = 0.682301827240309
$ws.Range("G5").Value = 1.564495524416255
$ws.Range("H5").Value = 1.001253100487087
$ws.Range("J5").Value = 0.6220977151477445
$ws.Range("G6").Value = 2.00711321810901
$ws.Range("H6").Value = 1.140267738152051
$ws.Range("J6").Value = 0.6064148509362324
$ws.Range("B7").Value = -2.518427412732279
$ws.Range("C7").Value = -1.13524266567742
$ws.Range("D7").Value = -1.186498229823468
$ws.Range("E7").Value = -1.564495524416255
$ws.Range("F7").Value = -2.00711321810901
$ws.Range("H7").Value = -1.221344854142138
$ws.Range("I7").Value = -2.349152721356929
$ws.Range("J7").Value = -1.989042107580228
$ws.Range("B8").Value = -1.416312553198512
$ws.Range("C8").Value = -1.030039409186893
$ws.Range("D8").Value = -1.105839399979577
$ws.Range("E8").Value = -1.001253100487087
$ws.Range("F8").Value = -1.140267738152051
$ws.Range("G8").Value = 1.221344854142138
$ws.Range("I8").Value = -0.7020039879691905
$ws.Range("J8").Value = -1.495991381269842
$ws.Range("G9").Value = 2.349152721356929
$ws.Range("H9").Value = 0.7020039879691905
$ws.Range("J9").Value = -0.3896871931125849
$ws.Range("B10").Value = -0.8625711600869045
$ws.Range("C10").Value = -0.1725125760959997
$ws.Range("D10").Value = -0.682301827240309
$ws.Range("E10").Value = -0.6220977151477445
$ws.Range("F10").Value = -0.6064148509362324
$ws.Range("G10").Value = 1.989042107580228
$ws.Range("H10").Value = 1.495991381269842
$ws.Range("I10").Value = 0.3896871931125849

Write-Output "Applied all changes"